# "updates from Oct 12th" - add two new donors (Embassy of Sweden, Latter-Day
# St Charities) and backfill the "Category" column for the donors that were
# added previously but never got a category value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill the "Category" column (column B) for existing donor rows 32-38 ---
$ws.Range("B32").Value = "Bilateral"
$ws.Range("B33").Value = "Foundation"
$ws.Range("B34").Value = "Private donors"
$ws.Range("B35").Value = "Foundation"
$ws.Range("B36").Value = "Combination"
$ws.Range("B37").Value = "Foundation"
$ws.Range("B38").Value = "Foundation"

# --- Append two brand-new donor rows ---
$ws.Range("A39").Value = "Embassy of Sweden"
$ws.Range("B39").Value = "Bilateral"

$ws.Range("A40").Value = "Latter-Day St Charities"
$ws.Range("B40").Value = "Foundation"

# --- The newly-entered cells (plus A39, which got caught up along with it)
#     carry a red font, matching the highlighting used to flag the edits ---
$redCells = @("B32","B33","B34","B35","B36","B37","B38","A39","B39","B40")
foreach ($addr in $redCells) {
    $ws.Range($addr).Font.Color = 255
}

# --- Leave the selection where the author left off, ready for the next entry ---
$ws.Range("B41").Select()
